# Update NCAA bracket results for final scores.
# Fills in Winner (C) / Loser (D) for every completed game, including the
# Round-of-32 / Sweet-16 / Elite-Eight / Final-Four / Title-game rows that
# were previously blank, and fixes two team-name spellings
# ("San Diego St." -> "San Diego St", "Utah St." -> "Utah State").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$results = @(
    @{ Row = 2;  Winner = "UConn";          Loser = "Stetson" }
    @{ Row = 3;  Winner = "Northwestern";   Loser = "FAU" }
    @{ Row = 4;  Winner = "San Diego St";   Loser = "UAB" }
    @{ Row = 5;  Winner = "Yale";           Loser = "Auburn" }
    @{ Row = 6;  Winner = "Duquesne";       Loser = "BYU" }
    @{ Row = 7;  Winner = "Illinois";       Loser = "Morehead St." }
    @{ Row = 8;  Winner = "Washington St."; Loser = "Drake" }
    @{ Row = 9;  Winner = "Iowa St.";       Loser = "South Dakota St." }
    @{ Row = 10; Winner = "UConn";          Loser = "Northwestern" }
    @{ Row = 11; Winner = "San Diego St";   Loser = "Yale" }
    @{ Row = 12; Winner = "Illinois";       Loser = "Duquesne" }
    @{ Row = 13; Winner = "Iowa St.";       Loser = "Washington St." }
    @{ Row = 14; Winner = "UConn";          Loser = "San Diego St" }
    @{ Row = 15; Winner = "Illinois";       Loser = "Iowa St." }
    @{ Row = 16; Winner = "UConn";          Loser = "Illinois" }
    @{ Row = 17; Winner = "Purdue";         Loser = "Grambling St." }
    @{ Row = 18; Winner = "Utah State";     Loser = "TCU" }
    @{ Row = 19; Winner = "Gonzaga";        Loser = "McNeese St." }
    @{ Row = 20; Winner = "Kansas";         Loser = "Samford" }
    @{ Row = 21; Winner = "Oregon";         Loser = "South Carolina" }
    @{ Row = 22; Winner = "Creighton";      Loser = "Akron" }
    @{ Row = 23; Winner = "Colorado St.";   Loser = "Texas" }
    @{ Row = 24; Winner = "Tennessee";      Loser = "St. Peter's" }
    @{ Row = 25; Winner = "Purdue";         Loser = "Utah State" }
    @{ Row = 26; Winner = "Gonzaga";        Loser = "Kansas" }
    @{ Row = 27; Winner = "Oregon";         Loser = "Creighton" }
    @{ Row = 28; Winner = "Tennessee";      Loser = "Colorado St." }
    @{ Row = 29; Winner = "Purdue";         Loser = "Gonzaga" }
    @{ Row = 30; Winner = "Tennessee";      Loser = "Creighton" }
    @{ Row = 31; Winner = "Purdue";         Loser = "Tennessee" }
    @{ Row = 32; Winner = "Houston";        Loser = "Longwood" }
    @{ Row = 33; Winner = "Texas A&M";      Loser = "Nebraska" }
    @{ Row = 34; Winner = "James Madison";  Loser = "Wisconsin" }
    @{ Row = 35; Winner = "Duke";           Loser = "Vermont" }
    @{ Row = 36; Winner = "NC State";       Loser = "Texas Tech" }
    @{ Row = 37; Winner = "Oakland";        Loser = "Kentucky" }
    @{ Row = 38; Winner = "Colorado";       Loser = "Florida" }
    @{ Row = 39; Winner = "Marquette";      Loser = "Western Kentucky" }
    @{ Row = 40; Winner = "Houston";        Loser = "Texas A&M" }
    @{ Row = 41; Winner = "Duke";           Loser = "James Madison" }
    @{ Row = 42; Winner = "NC State";       Loser = "Oakland" }
    @{ Row = 43; Winner = "Marquette";      Loser = "Colorado" }
    @{ Row = 44; Winner = "Duke";           Loser = "Houston" }
    @{ Row = 45; Winner = "NC State";       Loser = "Marquette" }
    @{ Row = 46; Winner = "NC State";       Loser = "Duke" }
    @{ Row = 47; Winner = "North Carolina"; Loser = "Wagner" }
    @{ Row = 48; Winner = "Michigan St.";   Loser = "Mississippi St." }
    @{ Row = 49; Winner = "Grand Canyon";   Loser = "St. Mary's" }
    @{ Row = 50; Winner = "Alabama";        Loser = "Charleston" }
    @{ Row = 51; Winner = "Clemson";        Loser = "New Mexico" }
    @{ Row = 52; Winner = "Baylor";         Loser = "Colgate" }
    @{ Row = 53; Winner = "Dayton";         Loser = "Nevada" }
    @{ Row = 54; Winner = "Arizona";        Loser = "Long Beach" }
    @{ Row = 55; Winner = "North Carolina"; Loser = "Michigan St." }
    @{ Row = 56; Winner = "Alabama";        Loser = "Grand Canyon" }
    @{ Row = 57; Winner = "Clemson";        Loser = "Baylor" }
    @{ Row = 58; Winner = "Arizona";        Loser = "Dayton" }
    @{ Row = 59; Winner = "Alabama";        Loser = "North Carolina" }
    @{ Row = 60; Winner = "Clemson";        Loser = "Arizona" }
    @{ Row = 61; Winner = "Alabama";        Loser = "Clemson" }
    @{ Row = 62; Winner = "UConn";          Loser = "Alabama" }
    @{ Row = 63; Winner = "Purdue";         Loser = "NC State" }
    @{ Row = 64; Winner = "UConn";          Loser = "Purdue" }
)

foreach ($r in $results) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Winner
    $ws.Cells.Item($r.Row, 4).Value = $r.Loser
}

# Match the author's last-saved selection.
$ws.Range("F24").Select() | Out-Null
